$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.452.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.800.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "39.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.288"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0668"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.059.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.804.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.417.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("E16").Value = "  -4.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.26%  "

$ws.Range("E20").Value = "  -4.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.82%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.45%  "

$ws.Range("E24").Value = "  -4.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.47%  "

$ws.Range("E28").Value = "  -1.77%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  -1.63%  "

$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("E32").Value = "  -3.60%  "

$ws.Range("E33").Value = "  -5.04%  "

$ws.Range("E34").Value = "  -0.99%  "

$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.306.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.03%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.44%  "

$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.69%  "

$ws.Range("E38").Value = "  -2.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.48%  "

$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "81.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "

$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.943"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.961.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.83%  "

$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.68%  "

$ws.Range("E51").Value = "  -0.54%  "
